$wb = $excel.ActiveWorkbook

# --- Sheet "Diff": add 5 new rows of package info, widen column B ---
$diff = $wb.Worksheets.Item("Diff")
$diff.Columns.Item(2).ColumnWidth = 44.139196

$diffRows = @(
    @("nuget", "Microsoft.AspNetCore.SpaServices.Extensions", "UNCHANGED", "8.0.23", "MIT", "8.0.23", "MIT", ""),
    @("nuget", "Microsoft.Extensions.FileProviders.Abstractions", "UNCHANGED", "8.0.0", "MIT", "8.0.0", "MIT", ""),
    @("nuget", "Microsoft.Extensions.FileProviders.Physical", "UNCHANGED", "8.0.0", "MIT", "8.0.0", "MIT", ""),
    @("nuget", "Microsoft.Extensions.FileSystemGlobbing", "UNCHANGED", "8.0.0", "MIT", "8.0.0", "MIT", ""),
    @("nuget", "Microsoft.Extensions.Primitives", "UNCHANGED", "8.0.0", "MIT", "8.0.0", "MIT", "")
)

for ($i = 0; $i -lt $diffRows.Count; $i++) {
    $rowIndex = $i + 2
    $rowData = $diffRows[$i]
    $arr = New-Object 'object[,]' 1,8
    for ($c = 0; $c -lt 8; $c++) {
        $arr[0,$c] = $rowData[$c]
    }
    $diff.Range($diff.Cells.Item($rowIndex, 1), $diff.Cells.Item($rowIndex, 8)).Value = $arr
}

# --- Sheet "CurrentDependencies": add same packages, rename headers, widen column B ---
$cur = $wb.Worksheets.Item("CurrentDependencies")
$cur.Columns.Item(2).ColumnWidth = 44.139196
$cur.Range("C1").Value = "Version"
$cur.Range("D1").Value = "License"

$curRows = @(
    @("nuget", "Microsoft.AspNetCore.SpaServices.Extensions", "8.0.23", "MIT", ""),
    @("nuget", "Microsoft.Extensions.FileProviders.Abstractions", "8.0.0", "MIT", ""),
    @("nuget", "Microsoft.Extensions.FileProviders.Physical", "8.0.0", "MIT", ""),
    @("nuget", "Microsoft.Extensions.FileSystemGlobbing", "8.0.0", "MIT", ""),
    @("nuget", "Microsoft.Extensions.Primitives", "8.0.0", "MIT", "")
)

for ($i = 0; $i -lt $curRows.Count; $i++) {
    $rowIndex = $i + 2
    $rowData = $curRows[$i]
    $arr = New-Object 'object[,]' 1,5
    for ($c = 0; $c -lt 5; $c++) {
        $arr[0,$c] = $rowData[$c]
    }
    $cur.Range($cur.Cells.Item($rowIndex, 1), $cur.Cells.Item($rowIndex, 5)).Value = $arr
}
